$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "243.98", "1.000").
# Setting .Value directly on a General-formatted cell would make Excel
# coerce them to actual numbers (losing exact text/trailing zeros), so
# each D cell is briefly marked as Text, written, then restored to the
# Normal style (matching the original formatting / style index).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.705.21'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.80%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.853.52'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.38%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9997'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.98'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6380'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '46.75'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.92%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07476'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.2979'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '24.30'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.29%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07643'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.51%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.851.14'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.038'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.81%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6873'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.45%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '83.69'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000009517'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +6.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.049'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '29.725.45'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.109.11'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '235.06'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.61%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.61'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9999'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.380'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.75%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.001'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.17'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1416'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.481'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '17.88'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.06272'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.492'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.273'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.142'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.084'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.174'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.853'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.31%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.7273'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.56%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.606'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.847'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01779'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.90%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.199.47'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9236'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.143'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.006.38'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.79'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '65.88'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.80%  '

# Rows 48-51 shift down by one: a new coin (BabyDogeCoin) is inserted at
# row 48 and the list is capped at 50 ranked rows, so RenderToken (which
# was previously row 51) drops off the bottom.

$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.00000000120'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.87%  '
$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4053'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.95%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.171'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.62%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05793'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.72%  '
